$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs", "Ntn1", "Unc5a", "ECs", 2, 0.6666666666666666, 0.8750386666666667, 2.625116, 0.05304058862308838, 0.05304058862308838, 3, 1, 2.322965, 6.968895, 0.2900587755589272, 0.2900587755589272, 2.032684196313333, 18.29415776682, 0.01538488819093778, 0.01538488819093778),
  @("ECs", "Ntn1", "Unc5a", "FAPs", 2, 0.6666666666666666, 0.8750386666666667, 2.625116, 0.05304058862308838, 0.05304058862308838, 3, 1, 1.626459, 4.879377, 0.2030890289077955, 0.2030890289077955, 1.423214514748, 12.808930632732, 0.01077196163616089, 0.01077196163616089),
  @("ECs", "Ntn1", "Unc5a", "M2", 2, 0.6666666666666666, 0.8750386666666667, 2.625116, 0.05304058862308838, 0.05304058862308838, 3, 1, 0.904191, 2.712573, 0.1129024907096758, 0.1129024907096758, 0.791202087052, 7.120818783468001, 0.00598841456425397, 0.005988414564253969),
  @("ECs", "Ntn1", "Unc5a", "sCs", 2, 0.6666666666666666, 0.8750386666666667, 2.625116, 0.05304058862308838, 0.05304058862308838, 3, 1, 3.154986, 9.464957999999999, 0.3939497048236015, 0.3939497048236015, 2.760734742792, 24.846612685128, 0.02089532423173575, 0.02089532423173574),
  @("FAPs", "Ntn1", "Unc5a", "ECs", 3, 1, 10.61942, 31.85826, 0.6436975977089742, 0.6436975977089742, 3, 1, 2.322965, 6.968895, 0.2900587755589272, 0.2900587755589272, 24.6685409803, 222.0168688227, 0.186710137021688, 0.186710137021688),
  @("FAPs", "Ntn1", "Unc5a", "FAPs", 3, 1, 10.61942, 31.85826, 0.6436975977089742, 0.6436975977089742, 3, 1, 1.626459, 4.879377, 0.2030890289077955, 0.2030890289077955, 17.27205123378, 155.44846110402, 0.1307279200289964, 0.1307279200289964),
  @("FAPs", "Ntn1", "Unc5a", "M2", 3, 1, 10.61942, 31.85826, 0.6436975977089742, 0.6436975977089742, 3, 1, 0.904191, 2.712573, 0.1129024907096758, 0.1129024907096758, 9.601983989219999, 86.41785590297999, 0.07267506204517807, 0.07267506204517805),
  @("FAPs", "Ntn1", "Unc5a", "sCs", 3, 1, 10.61942, 31.85826, 0.6436975977089742, 0.6436975977089742, 3, 1, 3.154986, 9.464957999999999, 0.3939497048236015, 0.3939497048236015, 33.50412142811999, 301.53709285308, 0.2535844786131118, 0.2535844786131117),
  @("M2", "Ntn1", "Unc5a", "ECs", 3, 1, 0.4029073333333333, 1.208722, 0.02442228319117198, 0.02442228319117198, 3, 1, 2.322965, 6.968895, 0.2900587755589272, 0.2900587755589272, 0.9359396335766665, 8.423456702189998, 0.007083897558784713, 0.007083897558784713),
  @("M2", "Ntn1", "Unc5a", "FAPs", 3, 1, 0.4029073333333333, 1.208722, 0.02442228319117198, 0.02442228319117198, 3, 1, 1.626459, 4.879377, 0.2030890289077955, 0.2030890289077955, 0.6553122584659998, 5.897810326193999, 0.004959897777006294, 0.004959897777006294),
  @("M2", "Ntn1", "Unc5a", "M2", 3, 1, 0.4029073333333333, 1.208722, 0.02442228319117198, 0.02442228319117198, 3, 1, 0.904191, 2.712573, 0.1129024907096758, 0.1129024907096758, 0.364305184634, 3.278746661705999, 0.002757336601100365, 0.002757336601100364),
  @("M2", "Ntn1", "Unc5a", "sCs", 3, 1, 0.4029073333333333, 1.208722, 0.02442228319117198, 0.02442228319117198, 3, 1, 3.154986, 9.464957999999999, 0.3939497048236015, 0.3939497048236015, 1.271166995964, 11.440502963676, 0.009621151254280605, 0.009621151254280605),
  @("sCs", "Ntn1", "Unc5a", "ECs", 3, 1, 4.600163333333334, 13.80049, 0.2788395304767656, 0.2788395304767655, 3, 1, 2.322965, 6.968895, 0.2900587755589272, 0.2900587755589272, 10.68601841761667, 96.17416575855, 0.08087985278751679, 0.08087985278751678),
  @("sCs", "Ntn1", "Unc5a", "FAPs", 3, 1, 4.600163333333334, 13.80049, 0.2788395304767656, 0.2788395304767655, 3, 1, 1.626459, 4.879377, 0.2030890289077955, 0.2030890289077955, 7.48197705497, 67.33779349472999, 0.05662924946563197, 0.05662924946563196),
  @("sCs", "Ntn1", "Unc5a", "M2", 3, 1, 4.600163333333334, 13.80049, 0.2788395304767656, 0.2788395304767655, 3, 1, 0.904191, 2.712573, 0.1129024907096758, 0.1129024907096758, 4.15942628453, 37.43483656077, 0.03148167749914338, 0.03148167749914337),
  @("sCs", "Ntn1", "Unc5a", "sCs", 3, 1, 4.600163333333334, 13.80049, 0.2788395304767656, 0.2788395304767655, 3, 1, 3.154986, 9.464957999999999, 0.3939497048236015, 0.3939497048236015, 14.51345091438, 130.62105822942, 0.1098487507244734, 0.1098487507244734),
)


for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
}

Write-Output "Updated $($data.Length) rows"
